$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the listed rows to reflect repulled data
$ws.Range("F2").Value = -3
$ws.Range("F6").Value = -4
$ws.Range("F8").Value = 0
$ws.Range("F10").Value = -5
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = 6
$ws.Range("F13").Value = 4
$ws.Range("F16").Value = -6
$ws.Range("F23").Value = -7
$ws.Range("F29").Value = -1
$ws.Range("F34").Value = -3
$ws.Range("F36").Value = 1
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 2
$ws.Range("F41").Value = 3
$ws.Range("F45").Value = 3
$ws.Range("F46").Value = -2
$ws.Range("F54").Value = 0
